$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add ignore Index option": renumber the index column (A) so it continues
# incrementing across the concatenated blocks instead of restarting at 0
# for each group (rows 5-10 become 3-8, continuing after rows 2-4 which are 0-2).
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
